$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append newly logged play-by-play yardage figures to the
# running space-separated lists stored in B2/C2 (row "R") and B3/C3 (row "P")
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$r = $wsYDS.Range("B2")
$r.Value2 = $r.Value2 + " 0 3 4 0 1 3 8 12 5 5 2 3 2 7 2 6 3 8 3 10 4 6 7 2 7 4 0 1 0 0 1 5 -2 3"

$r = $wsYDS.Range("C2")
$r.Value2 = $r.Value2 + " -2 0 2 1 21 5 5 8 6 6 6 12 4 6 3 6 6 5 8 4 2 11 17 2"

$r = $wsYDS.Range("B3")
$r.Value2 = $r.Value2 + " 10 11 11 9 13 -1 5 18 18 5 7 3 14 8 50 3 8 9 4 6 9 16 22"

$r = $wsYDS.Range("C3")
$r.Value2 = $r.Value2 + " 18 9 5 -3 24 10 0 1 3 5 10 6 2 17 4 8 2 5 6 6 3 3 7 4 6 12"

# ---------------------------------------------------------------------------
# OFF sheet: updated season totals after logging the week
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value2 = 413
$wsOFF.Range("D2").Value2 = 23
$wsOFF.Range("F2").Value2 = 96
$wsOFF.Range("G2").Value2 = 119
$wsOFF.Range("H2").Value2 = 5
$wsOFF.Range("I2").Value2 = 14
$wsOFF.Range("J2").Value2 = 52
$wsOFF.Range("L2").Value2 = 746
$wsOFF.Range("M2").Value2 = 484
$wsOFF.Range("O2").Value2 = 34
$wsOFF.Range("P2").Value2 = 20
$wsOFF.Range("Q2").Value2 = 1226
$wsOFF.Range("C3").Value2 = 460
$wsOFF.Range("D3").Value2 = 20
$wsOFF.Range("E3").Value2 = 65
$wsOFF.Range("F3").Value2 = 256
$wsOFF.Range("G3").Value2 = 119
$wsOFF.Range("H3").Value2 = 53
$wsOFF.Range("I3").Value2 = 112
$wsOFF.Range("J3").Value2 = 121
$wsOFF.Range("N3").Value2 = 29

# ---------------------------------------------------------------------------
# DEF sheet: updated season totals after logging the week
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value2 = 395
$wsDEF.Range("E2").Value2 = 6
$wsDEF.Range("F2").Value2 = 131
$wsDEF.Range("G2").Value2 = 117
$wsDEF.Range("J2").Value2 = 67
$wsDEF.Range("L2").Value2 = 643
$wsDEF.Range("M2").Value2 = 417
$wsDEF.Range("O2").Value2 = 56
$wsDEF.Range("P2").Value2 = 36
$wsDEF.Range("Q2").Value2 = 1144
$wsDEF.Range("C3").Value2 = 381
$wsDEF.Range("E3").Value2 = 54
$wsDEF.Range("F3").Value2 = 215
$wsDEF.Range("G3").Value2 = 105
$wsDEF.Range("I3").Value2 = 110
$wsDEF.Range("J3").Value2 = 112
$wsDEF.Range("N3").Value2 = 39

# ---------------------------------------------------------------------------
# ST sheet: updated season totals, plus appended weekly-log strings
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value2 = 197
$wsST.Range("D2").Value2 = 88
$wsST.Range("F2").Value2 = 245
$wsST.Range("G2").Value2 = 229
$wsST.Range("J2").Value2 = 93
$wsST.Range("K2").Value2 = 90
$wsST.Range("L2").Value2 = 43
$wsST.Range("M2").Value2 = 37
$wsST.Range("B3").Value2 = 138

$r = $wsST.Range("D3")
$r.Value2 = $r.Value2 + " 41 50"

$r = $wsST.Range("B4")
$r.Value2 = $r.Value2 + " 64 64 63"

$r = $wsST.Range("D4")
$r.Value2 = $r.Value2 + " 0 20"

$r = $wsST.Range("B5")
$r.Value2 = $r.Value2 + " 19 29 13"

$r = $wsST.Range("D5")
$r.Value2 = $r.Value2 + " 18 0"

$r = $wsST.Range("B6")
$r.Value2 = $r.Value2 + " 20"

# ---------------------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("C2").Value2 = 18
$wsTURNS.Range("D2").Value2 = 14
$wsTURNS.Range("E2").Value2 = 18
$wsTURNS.Range("D3").Value2 = 11
$wsTURNS.Range("E3").Value2 = 16

# ---------------------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value2 = 40
$wsPEN.Range("D2").Value2 = 15
$wsPEN.Range("D4").Value2 = 23
